# correct input data, edit text, draft new ncp pkgs
#
# The "iode_quality_flag" categorical variable (and its code lookup table)
# was removed from the workbook:
#   - ColumnHeadersToi!A17:C17 (the attribute row referencing iode_quality_flag)
#   - CategoricalVariables!A4:C8 (the code/definition rows for iode_quality_flag)
# Selections/active sheet are left the way the author left them when they
# finished (CategoricalVariables ends up the active tab).

$wb = $excel.ActiveWorkbook

# --- ColumnHeadersToi: remove the iode_quality_flag attribute row (row 17) ---
$wsToi = $wb.Worksheets.Item("ColumnHeadersToi")
$wsToi.Rows.Item(17).Delete()
$wsToi.Activate()
$wsToi.Range("A17:XFD17").Select()

# --- CategoricalVariables: remove the iode_quality_flag code table (rows 4-8) ---
$wsCat = $wb.Worksheets.Item("CategoricalVariables")
$wsCat.Range("A4:C8").EntireRow.Delete()
$wsCat.Activate()
$wsCat.Range("A4:D8").Select()
